# Scheduled-runner style update of market price / profit figures across the
# ALC, ARM, CRP, GSM, LTW and WVR leve-profit sheets. Columns H..N hold plain
# numeric (non-formula) values pulled from a market data source, so each
# changed cell is simply overwritten with its refreshed value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1577
$ws.Range("I19").Value = 1857.75
$ws.Range("J19").Value = 935.2857
$ws.Range("K19").Value = 1857.75
$ws.Range("L19").Value = 935.2857
$ws.Range("M19").Value = -1682.75
$ws.Range("N19").Value = -1285.2857

$ws.Range("H113").Value = 1793.5
$ws.Range("I113").Value = 1983.3334
$ws.Range("J113").Value = 1730.2222
$ws.Range("K113").Value = 1983.3334
$ws.Range("L113").Value = 1730.2222
$ws.Range("M113").Value = 1270.6666
$ws.Range("N113").Value = -8238.2222

$ws.Range("H129").Value = 1723.1052
$ws.Range("J129").Value = 1646.6111
$ws.Range("L129").Value = 4939.8333
$ws.Range("N129").Value = -14939.8333

$ws.Range("H132").Value = 2492.5435
$ws.Range("I132").Value = 2480.3953
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 7441.1859
$ws.Range("L132").Value = 8000.000100000001
$ws.Range("M132").Value = -4911.1859
$ws.Range("N132").Value = -13060.0001

$ws.Range("H137").Value = 23811180
$ws.Range("I137").Value = 1168.9688
$ws.Range("J137").Value = 100003220
$ws.Range("K137").Value = 3506.9064
$ws.Range("L137").Value = 300009660
$ws.Range("M137").Value = -956.9064000000003
$ws.Range("N137").Value = -300014760

$ws.Range("H138").Value = 2733.5
$ws.Range("I138").Value = 2807.75
$ws.Range("J138").Value = 2609.75
$ws.Range("K138").Value = 8423.25
$ws.Range("L138").Value = 7829.25
$ws.Range("M138").Value = -3283.25
$ws.Range("N138").Value = -18109.25

$ws.Range("H141").Value = 1149
$ws.Range("I141").Value = 929.53845
$ws.Range("J141").Value = 2100
$ws.Range("K141").Value = 2788.61535
$ws.Range("L141").Value = 6300
$ws.Range("M141").Value = 2391.38465
$ws.Range("N141").Value = -16660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7268.6987
$ws.Range("I32").Value = 5491.1514
$ws.Range("K32").Value = 5491.1514
$ws.Range("M32").Value = -5204.1514

$ws.Range("H74").Value = 2817.7693
$ws.Range("I74").Value = 1224.8
$ws.Range("J74").Value = 3813.375
$ws.Range("K74").Value = 1224.8
$ws.Range("L74").Value = 3813.375
$ws.Range("M74").Value = -350.8
$ws.Range("N74").Value = -5561.375

$ws.Range("H77").Value = 2817.7693
$ws.Range("I77").Value = 1224.8
$ws.Range("J77").Value = 3813.375
$ws.Range("K77").Value = 6124
$ws.Range("L77").Value = 19066.875
$ws.Range("M77").Value = -1756
$ws.Range("N77").Value = -27802.875

$ws.Range("H82").Value = 29800
$ws.Range("J82").Value = 29800
$ws.Range("L82").Value = 29800
$ws.Range("N82").Value = -30522

$ws.Range("H85").Value = 29800
$ws.Range("J85").Value = 29800
$ws.Range("L85").Value = 29800
$ws.Range("N85").Value = -32296

$ws.Range("H122").Value = 1165.1072
$ws.Range("I122").Value = 1206.4783
$ws.Range("J122").Value = 974.8
$ws.Range("K122").Value = 3619.4349
$ws.Range("L122").Value = 2924.4
$ws.Range("M122").Value = -1169.4349
$ws.Range("N122").Value = -7824.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1188.6666
$ws.Range("I7").Value = 3400
$ws.Range("J7").Value = 83
$ws.Range("K7").Value = 3400
$ws.Range("L7").Value = 83
$ws.Range("M7").Value = -3287
$ws.Range("N7").Value = -309

$ws.Range("H63").Value = 59000
$ws.Range("J63").Value = 59000
$ws.Range("L63").Value = 59000
$ws.Range("N63").Value = -60372

$ws.Range("H66").Value = 59000
$ws.Range("J66").Value = 59000
$ws.Range("L66").Value = 177000
$ws.Range("N66").Value = -183864

$ws.Range("H122").Value = 1336.5333
$ws.Range("I122").Value = 1029.3846
$ws.Range("J122").Value = 3333
$ws.Range("K122").Value = 3088.1538
$ws.Range("L122").Value = 9999
$ws.Range("M122").Value = -638.1538
$ws.Range("N122").Value = -14899

$ws.Range("H132").Value = 3491.85
$ws.Range("I132").Value = 3393.7273
$ws.Range("J132").Value = 3611.7778
$ws.Range("K132").Value = 10181.1819
$ws.Range("L132").Value = 10835.3334
$ws.Range("M132").Value = -7651.1819
$ws.Range("N132").Value = -15895.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 20000
$ws.Range("J68").Value = 20000
$ws.Range("L68").Value = 20000
$ws.Range("N68").Value = -21622

$ws.Range("H71").Value = 20000
$ws.Range("J71").Value = 20000
$ws.Range("L71").Value = 60000
$ws.Range("N71").Value = -68112

$ws.Range("H74").Value = 59610
$ws.Range("J74").Value = 59610
$ws.Range("L74").Value = 59610
$ws.Range("N74").Value = -61482

$ws.Range("H77").Value = 59610
$ws.Range("J77").Value = 59610
$ws.Range("L77").Value = 178830
$ws.Range("N77").Value = -188190

$ws.Range("H113").Value = 1200.5
$ws.Range("I113").Value = 882.25
$ws.Range("J113").Value = 1677.875
$ws.Range("K113").Value = 882.25
$ws.Range("L113").Value = 1677.875
$ws.Range("M113").Value = 1287.75
$ws.Range("N113").Value = -6017.875

$ws.Range("H132").Value = 2191.2144
$ws.Range("I132").Value = 1583.6666
$ws.Range("J132").Value = 2646.875
$ws.Range("K132").Value = 4750.9998
$ws.Range("L132").Value = 7940.625
$ws.Range("M132").Value = -2220.9998
$ws.Range("N132").Value = -13000.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1723.5
$ws.Range("I7").Value = 1631.3334
$ws.Range("K7").Value = 1631.3334
$ws.Range("M7").Value = -1519.3334

$ws.Range("H16").Value = 2493.1428
$ws.Range("I16").Value = 2827.4167
$ws.Range("J16").Value = 487.5
$ws.Range("K16").Value = 2827.4167
$ws.Range("L16").Value = 487.5
$ws.Range("M16").Value = -2657.4167
$ws.Range("N16").Value = -827.5

$ws.Range("H126").Value = 1723.5
$ws.Range("I126").Value = 1631.3334
$ws.Range("K126").Value = 4894.0002
$ws.Range("M126").Value = -2424.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1091.3715
$ws.Range("I122").Value = 1099.2727
$ws.Range("J122").Value = 1078
$ws.Range("K122").Value = 3297.8181
$ws.Range("L122").Value = 3234
$ws.Range("M122").Value = -847.8181
$ws.Range("N122").Value = -8134
